$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to make room for the new year (2022) entry.
# Excel shifts all existing data (rows 2-23) down to rows 3-24, which also
# naturally carries the "Fonte" source-link annotation from O2 down to O3.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the 2022 data point (only January/"B" column
# has a value so far; the remaining months are not yet reported).
$ws.Range("A2").Value = 2022
$ws.Range("B2").Value = 0.73229999999999995

# Remove the placeholder cells left behind by the insert for columns C..O
# on row 2 so that row only contains the A2/B2 cells, matching the source
# data (no empty cells recorded for months that have no data yet).
$ws.Range("C2:O2").Clear()

# Correct the December (M) value for year 2021, which now lives on row 3
# after the insert (was previously 0.77, corrected to 0.7691).
$ws.Range("M3").Value = 0.76910000000000001
